$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.548.80"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.811.98"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "663.12"
$ws.Range("E5").Value = "  +5.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.40"
$ws.Range("E6").Value = "  +2.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.812.26"
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.462"
$ws.Range("E11").Value = "  +2.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.98"
$ws.Range("E12").Value = "  +5.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000244"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.86"
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.450.51"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.813.57"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.429.34"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.82"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.19"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.40"
$ws.Range("E21").Value = "  +8.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "476.47"
$ws.Range("E22").Value = "  +1.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.715"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000147"
$ws.Range("E24").Value = "  -2.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.80"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.21"
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.40"
$ws.Range("E27").Value = "  +4.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.14"
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.959.51"
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.83"
$ws.Range("E31").Value = "  +5.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.33"
$ws.Range("E32").Value = "  +3.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.43"
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.35"
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("E35").Value = "  +17.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.764.83"
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "9.08"
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.36"
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.94"
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.973"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.09"
$ws.Range("E45").Value = "  +7.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.19"
$ws.Range("E46").Value = "  +5.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "158.65"
$ws.Range("E47").Value = "  +3.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "47.91"
$ws.Range("E48").Value = "  +2.29%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.43"
$ws.Range("E49").Value = "  +3.21%  "
$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.302"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.52"
$ws.Range("E51").Value = "  +1.04%  "
